# Update the "dSF" column (F) values on the active sheet to reflect
# repulled / recalculated data, per commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = 3
    4  = 3
    5  = 1
    6  = 6
    7  = 1
    8  = 1
    12 = -3
    15 = -1
    17 = -5
    19 = -3
    20 = 2
    21 = -4
    22 = 0
    25 = -1
    26 = -3
    28 = 12
    29 = 3
    30 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
